$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6442526578903198
$ws.Range("B1").Value = 2.268299341201782
$ws.Range("C1").Value = 6.268117904663086
$ws.Range("D1").Value = 1.7615567445755
$ws.Range("E1").Value = 1.04323410987854
